# Re-pull / push updated "dSF" (column F) data into the save-data sheet.
# Only column F values change for a subset of rows (mean calculation
# re-run on freshly re-pulled data); every other column is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    4  = 0
    10 = -1
    12 = 2
    15 = 1
    19 = -1
    21 = -1
    26 = -1
    37 = 1
    39 = 3
    42 = 2
    45 = -2
    47 = 4
    48 = -3
    52 = -3
    56 = 5
    57 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
